# "Short slides for Colin"
#
# 1) The footer "datetimeFigureOut" field cached on the slide master and
#    every slide layout is bumped from 10/28/2016 to 1/23/2017 (PowerPoint
#    re-caches this field's displayed text whenever the deck is saved).
# 2) Slide 2's "OAuth," bullet loses its trailing comma ("OAuth,"->"OAuth").
# 3) Embedded OLE worksheet objects (slides 5-8) keep their content but
#    PowerPoint renumbers their legacy VML shape ids (spid) on save; that
#    id isn't part of the automation object model, so it isn't touched here.

$p = $ppt.ActivePresentation

$newDate = "1/23/2017"

function Set-DatePlaceholderText($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide master date placeholder.
Set-DatePlaceholderText $p.SlideMaster.Shapes

# Every slide layout's date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Set-DatePlaceholderText $layouts.Item($li).Shapes
}

# Slide 2: "OAuth," -> "OAuth" (drop the trailing comma).
$slide2 = $p.Slides.Item(2)
$contentShape = $slide2.Shapes.Item(2)
$tr = $contentShape.TextFrame.TextRange
$paraCount = $tr.Paragraphs().Count
for ($j = 1; $j -le $paraCount; $j++) {
    $para = $tr.Paragraphs($j)
    if ($para.Text -eq "OAuth,") {
        $lastChar = $para.Characters($para.Length, 1)
        $lastChar.Delete()
    }
}
